# Add four new translation rows (Save changes / order confirmation / payment
# failure / order number strings) to the bottom of the language sheet,
# matching the style already used for the "label" column (column B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 160; B = "Save changes"; C = "שמור    וצא" },
    @{ Row = 161; B = "Your order has been accepted payment succed"; C = "הזמנתך התקבלה בהצלחה התשלום בוצע" },
    @{ Row = 162; B = "Sorry yout payment was not accepted. Order not palced"; C = "מצטערים התשלום לא התקבל ולפיכך ההזמנה לא  תבוצע" },
    @{ Row = 163; B = "order Number"; C = "קוד ההזמנה" }
)

# Existing row 159 column B carries the style (font/alignment) that every
# new label cell should reuse - copy its formatting instead of re-building
# the font manually so no new style entries are introduced.
$formatSource = $ws.Range("B159")

foreach ($r in $rows) {
    $bCell = $ws.Cells.Item($r.Row, 2)
    $cCell = $ws.Cells.Item($r.Row, 3)

    $bCell.Value = $r.B
    $cCell.Value = $r.C

    $formatSource.Copy()
    $bCell.PasteSpecial(-4122)
}

[void]$ws.Range("I159").Select()
$ws.Application.CutCopyMode = $false

Write-Output "Added rows 160-163 with new translation strings"
